$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("succinate_hopa")
$ws.Range("F2").Value = "SUCCSEMIALDDEHYDROG-RXN_reverse: [A0PN13, A0PST9, A0QMB9, A0R4Q0, A1KF54, A1KJE8, A5TYV9, A5U390, B1XMM6, O32507, O69497, O85997, P25526, P38067, P55653, P94428, P9WNX6, P9WNX7, P9WNX8, Q0K2K1, Q4KKA2, Q55585, Q73TP5, Q7TZP3, Q9JTN7, Q9RBF6] | R00714_reverse: [] | 1.2.1.79 ()_6: [] | 1.2.1.24_53: [] | 1.2.1.16_25: [] | "
$ws.Range("G2").Value = "RXN-11727: [] | R09545_reverse: [] | 1.1.1.308_0: [] | "
$ws.Range("H2").Value = "R07650_reverse: [] | 4.1.1.86_0: [] | 4.1.1.86-RXN: [D0S1M8, P71362, Q1V6P8, Q43908] | R06614_reverse: [] | R00656_reverse: [] | METHIONINE-DECARBOXYLASE-RXN: [] | 4.1.1.57_2: [] | 4.1.1.14_03: [] | RXN-7704: [] | 4.1.1.72_19: [] | 4.1.1.43_05: [] | 4.1.1.1_19: [] | 4.1.1.86_4: [] | 4.1.1.72_33: [] | 4.1.1.43_11: [] | 4.1.1.1_38: [] | 4.1.1.72_34: [] | 4.1.1.1_39: [] | RXN-14985: [P26263] | 4.1.1.72_29: [] | 4.1.1.1_34: [] | 4.1.1.7_151: [] | 4.1.1.72_37: [] | 4.1.1.1_41: [] | R03894_reverse: [] | 4.1.1.7_062: [] | 4.1.1.72_04: [] | 4.1.1.72_03: [] | 4.1.1.72-RXN: [P26263, Q06408, Q07471] | 4.1.1.1_04: [] | VALINE-DECARBOXYLASE-RXN: [] | R01437_reverse: [] | 4.1.1.57_1: [] | 4.1.1.14_01: [] | RXN-7692: [A2XFI3, A2Y5L9, A2YQ76, O42873, O82647, P16467, P26263, P28516, P33149, P33287, P34734, P51844, P51845, P51846, P51850, P51851, P83779, P87208, Q05326, Q05327, Q07471, Q09737, Q0CNV1, Q0D3D2, Q0DHF6, Q10MW3, Q12629, Q2UKV4, Q4WXX9, Q684J7, Q6FJA3, Q92345, Q9FFT4, Q9M039, Q9M040, Q9P7P6] | 4.1.1.7_096: [] | 4.1.1.72_17: [] | 4.1.1.1_17: [] | RXN-7643: [P16467, P26263] | 4.1.1.7_145: [] | 4.1.1.72_32: [] | 4.1.1.1_37: [] | 2.2.1.6_12: [] | "
$ws.Range("F3").Value = "RXN-11727: [] | R09545_reverse: [] | 1.1.1.308_0: [] | "
$ws.Range("G3").Value = "RXN-13089_reverse: [H8ZPX2] | "
$ws.Range("H3").Value = "R07650_reverse: [] | 4.1.1.86_0: [] | 4.1.1.86-RXN: [D0S1M8, P71362, Q1V6P8, Q43908] | R06614_reverse: [] | R00656_reverse: [] | METHIONINE-DECARBOXYLASE-RXN: [] | 4.1.1.57_2: [] | 4.1.1.14_03: [] | RXN-7704: [] | 4.1.1.72_19: [] | 4.1.1.43_05: [] | 4.1.1.1_19: [] | 4.1.1.86_4: [] | 4.1.1.72_33: [] | 4.1.1.43_11: [] | 4.1.1.1_38: [] | 4.1.1.72_34: [] | 4.1.1.1_39: [] | RXN-14985: [P26263] | 4.1.1.72_29: [] | 4.1.1.1_34: [] | 4.1.1.7_151: [] | 4.1.1.72_37: [] | 4.1.1.1_41: [] | R03894_reverse: [] | 4.1.1.7_062: [] | 4.1.1.72_04: [] | 4.1.1.72_03: [] | 4.1.1.72-RXN: [P26263, Q06408, Q07471] | 4.1.1.1_04: [] | VALINE-DECARBOXYLASE-RXN: [] | R01437_reverse: [] | 4.1.1.57_1: [] | 4.1.1.14_01: [] | RXN-7692: [A2XFI3, A2Y5L9, A2YQ76, O42873, O82647, P16467, P26263, P28516, P33149, P33287, P34734, P51844, P51845, P51846, P51850, P51851, P83779, P87208, Q05326, Q05327, Q07471, Q09737, Q0CNV1, Q0D3D2, Q0DHF6, Q10MW3, Q12629, Q2UKV4, Q4WXX9, Q684J7, Q6FJA3, Q92345, Q9FFT4, Q9M039, Q9M040, Q9P7P6] | 4.1.1.7_096: [] | 4.1.1.72_17: [] | 4.1.1.1_17: [] | RXN-7643: [P16467, P26263] | 4.1.1.7_145: [] | 4.1.1.72_32: [] | 4.1.1.1_37: [] | 2.2.1.6_12: [] | "

$ws = $wb.Worksheets.Item("oxaloacetate_hopa")
$ws.Range("F2").Value = "RXN-11727: [] | R09545_reverse: [] | 1.1.1.308_0: [] | "
$ws.Range("G2").Value = "R00930_reverse: [] | 2.1.3.1_7: [] | 2.1.3.1-RXN: [P02904, Q70AC7, Q8GBW6] | "
$ws.Range("H2").Value = "R07219_reverse: [] | 1.17.7.4_1: [] | "
$ws.Range("F3").Value = "RXN-14468_reverse: [] | 2.1.1.281 ()_1: [] | RXN-11456_reverse: [D5FKJ3] | R09939: [] | "
$ws.Range("G3").Value = "RXN-16543_reverse: [] | "
$ws.Range("H3").Value = "RXN-14984: [] | RXN-14984_reverse: [] | "
$ws.Range("F4").Value = "RXN-11727: [] | R09545_reverse: [] | 1.1.1.308_0: [] | "
$ws.Range("G4").Value = "RXN-14468_reverse: [] | 2.1.1.281 ()_1: [] | RXN-11456_reverse: [D5FKJ3] | R09939: [] | "
$ws.Range("H4").Value = "RXN-14984: [] | RXN-14984_reverse: [] | "

$ws = $wb.Worksheets.Item("pyruvate_hopa")
$ws.Range("F2").Value = "RXN-12074: [Q88JX9] | R00008_reverse: [] | 4.1.3.17_7: [] | 4.1.3.17-RXN: [] | "
$ws.Range("G2").Value = "R02272: [] | GSAAMINOTRANS-RXN: [A0AJ02, A0AJD9, A0B885, A0JYM2, A0K5D7, A0KNY9, A0KZS6, A0L3S9, A0LIP3, A0LRE5, A0LXY3, A0PLS0, A0Q2A8, A0Q628, A0QLG2, A0QR33, A0R9G6, A0RJ78, A0RQ51, A0RXB3, A1A7K0, A1ASE9, A1AVQ1, A1BJG8, A1JJQ1, A1KB59, A1KG00, A1KS34, A1R898, A1RMG7, A1RQR5, A1S3U0, A1SE06, A1ST92, A1T3F2, A1TKK8, A1TYG5, A1UAR1, A1V1L0, A1V9X6, A1VSM6, A1VZJ6, A1WBR8, A1WMA0, A1WYL3, A2BL27, A2BPW6, A2BVE5, A2C0U2, A2C7I7, A2S9D9, A2SKQ7, A2SQU3, A3CL81, A3CU65, A3D1R5, A3DIF3, A3M7I8, A3MMQ8, A3MWW7, A3N2K3, A3NCF3, A3NY82, A3PBK8, A3PUB7, A3QBN5, A4FPX3, A4FZY1, A4G8N0, A4IKL1, A4IRG2, A4IXN8, A4J6H0, A4JC86, A4QB78, A4SGT2, A4SJ79, A4SVE6, A4T358, A4TPW6, A4VQY0, A4W6Q1, A4WN33, A4XC73, A4XK09, A4XYV6, A4Y4G4, A4YD91, A5CNI7, A5CXH8, A5D3M2, A5EXX1, A5F945, A5FIM3, A5G9C0, A5GMT2, A5GUJ2, A5IC29, A5ITJ2, A5IU30, A5N5Y1, A5TZQ4, A5UML0, A5UU40, A5VM65, A5W9H0, A5WC94, A6GYW4, A6LKX4, A6LSY8, A6Q2X6, A6Q9P9, A6QHK1, A6QI46, A6T2A2, A6T4V8, A6TJD8, A6U2D5, A6U2W8, A6URL2, A6UWB3, A6V0D2, A6VIL1, A6VQ66, A6VUX8, A6W5M6, A6WKL6, A7FM09, A7GL00, A7GTE2, A7GXK4, A7H3N1, A7HIX3, A7I252, A7I7P2, A7MGR5, A7MUU9, A7NCX8, A7NKV1, A7X388, A7X3X9, A7Z2N8, A7Z7A2, A7ZCH5, A7ZHP6, A7ZWA2, A8AAB3, A8ALD5, A8ETJ2, A8FB92, A8FFU9, A8FLR2, A8FS76, A8G3J9, A8G9U7, A8H164, A8LCD0, A8M6S5, A8M919, A8MGY8, A8YY36, A8Z2I8, A8ZZU9, A9A1A0, A9A875, A9AEU0, A9B550, A9BEA5, A9BY74, A9GDD3, A9I1R3, A9KBA0, A9KP86, A9KSC7, A9L5J5, A9M1G6, A9MPK7, A9N0P9, A9NAX3, A9R1E6, A9VIT5, A9VSS7, A9WIS7, A9WPI0, B0B9W0, B0BBJ0, B0BRF0, B0CC57, B0JPW6, B0KJV9, B0R7L9, B0RBL3, B0RP74, B0SF64, B0SNK3, B0TFV0, B0TIQ0, B0U0T6, B0U328, B0UST3, B0VCN6, B0VU92, B1HUT1, B1HVD3, B1I4L1, B1IQI6, B1J1Y6, B1JK22, B1JX81, B1KHF6, B1LGV7, B1MHH1, B1VEK3, B1VTD1, B1X023, B1XD24, B1XIT5, B1XT79, B1Y8G2, B1YBL3, B1YJU8, B1YK58, B1YUJ7, B1ZVF7, B2A1H4, B2FT35, B2G9H2, B2GL27, B2HRL5, B2HW19, B2I6C2, B2J7M9, B2JED5, B2K547, B2SGG7, B2SKS0, B2T6C2, B2TPD4, B2U2Z9, B2U6Q7, B2UNE3, B2USD7, B2UX41, B2V5U0, B2VE25, B3DY05, B3E9X0, B3EI07, B3EKJ7, B3GYN8, B3PL68, B3QRD2, B3QSA6, B3R312, B4EBT3, B4EUE1, B4RNX2, B4S0D0, B4S3Q6, B4SGW1, B4SI23, B4SUY4, B4TK30, B4TXQ6, B4U8M6, B4UKE6, B5BL82, B5EFG4, B5ERN9, B5F8R5, B5FAL8, B5FJ01, B5R3G6, B5RHE0, B5Y1L5, B5YHH6, B5Z0D4, B5ZA72, B6EL04, B6HZD0, B6J3H2, B6J495, B6JKN5, B7GH35, B7GIK0, B7GYZ5, B7H9S6, B7HE93, B7HQM1, B7HU66, B7I4N4, B7IHH3, B7IIX2, B7IW22, B7JBI1, B7JNG6, B7JQ54, B7K2I1, B7KA18, B7LGL6, B7LWB5, B7M195, B7MBD7, B7MP16, B7N823, B7NIB7, B7UIK1, B7V9D9, B7VJJ5, B8CJG4, B8DFN7, B8DHJ5, B8DN02, B8EBU3, B8F7H5, B8FBA1, B8G2L6, B8G822, B8GJJ2, B8GLE4, B8HBR2, B8HYK1, B8I0R1, B8J3A3, B8JCU6, B8ZT67, B9DMX3, B9DND8, B9DZG0, B9E740, B9E7Y3, B9IZ36, B9J3N2, B9KFS8, B9L0Q2, B9LA83, B9LKS0, B9LTZ9, B9LZL6, B9MF68, B9MRJ7, C0Q5R5, C0QCR9, C0QQ98, C0Z671, C0ZAH4, C0ZUX9, C1A7K7, C1AKK4, C1AZ30, C1CWI8, C1D9B6, C1DMY5, C1DUY4, C1ETQ7, C1EWG6, C1F911, C1KVK0, C1KVY4, C3K2K3, C3L6Z2, C3LHA0, C3LSN1, C3MJ22, C3MYE0, C3MZR9, C3N842, C3NF69, C3P9E6, C3PCZ3, C3PKI4, C4K6Y5, C4KJ00, C4L367, C4L4J9, C4LAK4, C4ZRP7, C5BAP9, C6AR33, C6DC32, E1W874, O26330, O29027, O66998, O74038, O74061, O84212, P0C1P8, P0C2D9, P0CL07, P18492, P21267, P23893, P24630, P30949, P31593, P45621, P46716, P48247, P56115, P63507, P63508, P71084, P99096, P9WMN8, P9WMN9, Q01YQ2, Q02SE5, Q04NV4, Q04X52, Q06741, Q06774, Q07YU5, Q0AAH7, Q0AFV1, Q0AZ36, Q0BHJ2, Q0BLC7, Q0HG53, Q0HSE6, Q0I3R7, Q0I8G1, Q0KDN9, Q0RS03, Q0SF42, Q0ST19, Q0T852, Q0TLH7, Q0TQG7, Q0VSP6, Q0W5T3, Q110Z9, Q11WK1, Q12EF7, Q12KC8, Q12WM7, Q13V40, Q14HS1, Q15YL3, Q17X40, Q18ES5, Q1AUK6, Q1BE73, Q1BYB1, Q1C3X1, Q1CLU7, Q1CUJ7, Q1DFA7, Q1GXW0, Q1I4H5, Q1IHV2, Q1IWZ8, Q1LQM3, Q1MPW7, Q1QEE8, Q1QSB0, Q1RG34, Q21MJ0, Q21YQ0, Q24VD4, Q254Z2, Q2A2U8, Q2FFN1, Q2FG69, Q2FTK5, Q2FXR4, Q2G283, Q2IHP8, Q2JFQ1, Q2JMP7, Q2JS70, Q2KUS4, Q2LVW6, Q2NF42, Q2NVQ0, Q2P6E4, Q2RJ31, Q2S1S3, Q2S8X9, Q2SYB2, Q2Y5Q5, Q2YTC2, Q2YU22, Q30RJ4, Q30WH2, Q31C50, Q31FB1, Q31QJ2, Q325Y5, Q32JV4, Q39566, Q39IQ4, Q39QA6, Q3A7W5, Q3ACS9, Q3ALU9, Q3AP59, Q3AWP4, Q3B1A1, Q3BPP7, Q3IHS1, Q3IT20, Q3JDH7, Q3JPN1, Q3K6C0, Q3KMF2, Q3M3B9, Q3SFU6, Q3Z5K3, Q40147, Q42522, Q46CH1, Q46GT9, Q46XZ7, Q478V1, Q47LB3, Q47V96, Q48DP1, Q49Y99, Q49YQ9, Q4FVC4, Q4JAM7, Q4JSY1, Q4K5I6, Q4L706, Q4L7G9, Q4UYA7, Q4ZNA0, Q55665, Q57T53, Q58020, Q5E2W6, Q5FAH9, Q5H3I5, Q5HER0, Q5HFA5, Q5HN71, Q5HNN6, Q5HUU3, Q5KWK9, Q5L2S4, Q5L5P0, Q5NGB9, Q5P6J6, Q5PD43, Q5QVP9, Q5SJS4, Q5UZ90, Q5WEP8, Q5WIE9, Q5WWG1, Q5X529, Q5YP79, Q5ZVA6, Q60CV0, Q62HV8, Q633Y3, Q63GB4, Q63RP8, Q65GK4, Q65M82, Q65U43, Q66EF1, Q67KM3, Q6AB08, Q6AHE8, Q6AQ32, Q6D1Z0, Q6FCY1, Q6G870, Q6G8Q8, Q6GFJ3, Q6GG38, Q6HD65, Q6HNS8, Q6L2G9, Q6LUS3, Q6M0P5, Q6MAC7, Q6MHT9, Q6NJJ2, Q6YZE2, Q71YY2, Q71ZB5, Q725I1, Q72K83, Q72ZW5, Q73DX4, Q73SQ3, Q74GA9, Q75G04, Q7A4T5, Q7M847, Q7MHY9, Q7N845, Q7NPI4, Q7P1Z5, Q7U598, Q7UPM9, Q7V2J3, Q7V677, Q7VDA1, Q7VHK3, Q7W050, Q7W3U1, Q7WF71, Q817R3, Q81I85, Q81LD0, Q81YV0, Q821C1, Q822Q0, Q82E21, Q82UQ8, Q83AK3, Q83FJ5, Q83H98, Q85WB7, Q87BW3, Q87LY3, Q87VY5, Q88DP0, Q897K4, Q8CNZ1, Q8CRW7, Q8CV56, Q8CZE3, Q8D3C8, Q8DBX8, Q8DLK8, Q8EHC8, Q8EY44, Q8FL16, Q8FSD4, Q8KAQ7, Q8NT73, Q8NVU6, Q8NW75, Q8P5R4, Q8PH40, Q8PW58, Q8RFY7, Q8TT57, Q8TYL6, Q8X4V5, Q8Y1M4, Q8Y6J9, Q8Y6X8, Q8YS26, Q8Z9B4, Q8ZBL9, Q8ZYW1, Q92AX5, Q92BG1, Q976H2, Q97B25, Q97MU2, Q980U5, Q99T15, Q9CNG9, Q9F2S0, Q9HKM6, Q9HMY8, Q9JRW9, Q9JW10, Q9JXW0, Q9K8G3, Q9KEB0, Q9KU97, Q9PB43, Q9PKI3, Q9PP70, Q9RWW0, Q9Y9I9, Q9ZMD0] | 5.4.3.8_1: [] | "
$ws.Range("H2").Value = "RXN-17561_reverse: [] | "
$ws.Range("F3").Value = "RXN-12074: [Q88JX9] | R00008_reverse: [] | 4.1.3.17_7: [] | 4.1.3.17-RXN: [] | "
$ws.Range("G3").Value = "R02272: [] | GSAAMINOTRANS-RXN: [A0AJ02, A0AJD9, A0B885, A0JYM2, A0K5D7, A0KNY9, A0KZS6, A0L3S9, A0LIP3, A0LRE5, A0LXY3, A0PLS0, A0Q2A8, A0Q628, A0QLG2, A0QR33, A0R9G6, A0RJ78, A0RQ51, A0RXB3, A1A7K0, A1ASE9, A1AVQ1, A1BJG8, A1JJQ1, A1KB59, A1KG00, A1KS34, A1R898, A1RMG7, A1RQR5, A1S3U0, A1SE06, A1ST92, A1T3F2, A1TKK8, A1TYG5, A1UAR1, A1V1L0, A1V9X6, A1VSM6, A1VZJ6, A1WBR8, A1WMA0, A1WYL3, A2BL27, A2BPW6, A2BVE5, A2C0U2, A2C7I7, A2S9D9, A2SKQ7, A2SQU3, A3CL81, A3CU65, A3D1R5, A3DIF3, A3M7I8, A3MMQ8, A3MWW7, A3N2K3, A3NCF3, A3NY82, A3PBK8, A3PUB7, A3QBN5, A4FPX3, A4FZY1, A4G8N0, A4IKL1, A4IRG2, A4IXN8, A4J6H0, A4JC86, A4QB78, A4SGT2, A4SJ79, A4SVE6, A4T358, A4TPW6, A4VQY0, A4W6Q1, A4WN33, A4XC73, A4XK09, A4XYV6, A4Y4G4, A4YD91, A5CNI7, A5CXH8, A5D3M2, A5EXX1, A5F945, A5FIM3, A5G9C0, A5GMT2, A5GUJ2, A5IC29, A5ITJ2, A5IU30, A5N5Y1, A5TZQ4, A5UML0, A5UU40, A5VM65, A5W9H0, A5WC94, A6GYW4, A6LKX4, A6LSY8, A6Q2X6, A6Q9P9, A6QHK1, A6QI46, A6T2A2, A6T4V8, A6TJD8, A6U2D5, A6U2W8, A6URL2, A6UWB3, A6V0D2, A6VIL1, A6VQ66, A6VUX8, A6W5M6, A6WKL6, A7FM09, A7GL00, A7GTE2, A7GXK4, A7H3N1, A7HIX3, A7I252, A7I7P2, A7MGR5, A7MUU9, A7NCX8, A7NKV1, A7X388, A7X3X9, A7Z2N8, A7Z7A2, A7ZCH5, A7ZHP6, A7ZWA2, A8AAB3, A8ALD5, A8ETJ2, A8FB92, A8FFU9, A8FLR2, A8FS76, A8G3J9, A8G9U7, A8H164, A8LCD0, A8M6S5, A8M919, A8MGY8, A8YY36, A8Z2I8, A8ZZU9, A9A1A0, A9A875, A9AEU0, A9B550, A9BEA5, A9BY74, A9GDD3, A9I1R3, A9KBA0, A9KP86, A9KSC7, A9L5J5, A9M1G6, A9MPK7, A9N0P9, A9NAX3, A9R1E6, A9VIT5, A9VSS7, A9WIS7, A9WPI0, B0B9W0, B0BBJ0, B0BRF0, B0CC57, B0JPW6, B0KJV9, B0R7L9, B0RBL3, B0RP74, B0SF64, B0SNK3, B0TFV0, B0TIQ0, B0U0T6, B0U328, B0UST3, B0VCN6, B0VU92, B1HUT1, B1HVD3, B1I4L1, B1IQI6, B1J1Y6, B1JK22, B1JX81, B1KHF6, B1LGV7, B1MHH1, B1VEK3, B1VTD1, B1X023, B1XD24, B1XIT5, B1XT79, B1Y8G2, B1YBL3, B1YJU8, B1YK58, B1YUJ7, B1ZVF7, B2A1H4, B2FT35, B2G9H2, B2GL27, B2HRL5, B2HW19, B2I6C2, B2J7M9, B2JED5, B2K547, B2SGG7, B2SKS0, B2T6C2, B2TPD4, B2U2Z9, B2U6Q7, B2UNE3, B2USD7, B2UX41, B2V5U0, B2VE25, B3DY05, B3E9X0, B3EI07, B3EKJ7, B3GYN8, B3PL68, B3QRD2, B3QSA6, B3R312, B4EBT3, B4EUE1, B4RNX2, B4S0D0, B4S3Q6, B4SGW1, B4SI23, B4SUY4, B4TK30, B4TXQ6, B4U8M6, B4UKE6, B5BL82, B5EFG4, B5ERN9, B5F8R5, B5FAL8, B5FJ01, B5R3G6, B5RHE0, B5Y1L5, B5YHH6, B5Z0D4, B5ZA72, B6EL04, B6HZD0, B6J3H2, B6J495, B6JKN5, B7GH35, B7GIK0, B7GYZ5, B7H9S6, B7HE93, B7HQM1, B7HU66, B7I4N4, B7IHH3, B7IIX2, B7IW22, B7JBI1, B7JNG6, B7JQ54, B7K2I1, B7KA18, B7LGL6, B7LWB5, B7M195, B7MBD7, B7MP16, B7N823, B7NIB7, B7UIK1, B7V9D9, B7VJJ5, B8CJG4, B8DFN7, B8DHJ5, B8DN02, B8EBU3, B8F7H5, B8FBA1, B8G2L6, B8G822, B8GJJ2, B8GLE4, B8HBR2, B8HYK1, B8I0R1, B8J3A3, B8JCU6, B8ZT67, B9DMX3, B9DND8, B9DZG0, B9E740, B9E7Y3, B9IZ36, B9J3N2, B9KFS8, B9L0Q2, B9LA83, B9LKS0, B9LTZ9, B9LZL6, B9MF68, B9MRJ7, C0Q5R5, C0QCR9, C0QQ98, C0Z671, C0ZAH4, C0ZUX9, C1A7K7, C1AKK4, C1AZ30, C1CWI8, C1D9B6, C1DMY5, C1DUY4, C1ETQ7, C1EWG6, C1F911, C1KVK0, C1KVY4, C3K2K3, C3L6Z2, C3LHA0, C3LSN1, C3MJ22, C3MYE0, C3MZR9, C3N842, C3NF69, C3P9E6, C3PCZ3, C3PKI4, C4K6Y5, C4KJ00, C4L367, C4L4J9, C4LAK4, C4ZRP7, C5BAP9, C6AR33, C6DC32, E1W874, O26330, O29027, O66998, O74038, O74061, O84212, P0C1P8, P0C2D9, P0CL07, P18492, P21267, P23893, P24630, P30949, P31593, P45621, P46716, P48247, P56115, P63507, P63508, P71084, P99096, P9WMN8, P9WMN9, Q01YQ2, Q02SE5, Q04NV4, Q04X52, Q06741, Q06774, Q07YU5, Q0AAH7, Q0AFV1, Q0AZ36, Q0BHJ2, Q0BLC7, Q0HG53, Q0HSE6, Q0I3R7, Q0I8G1, Q0KDN9, Q0RS03, Q0SF42, Q0ST19, Q0T852, Q0TLH7, Q0TQG7, Q0VSP6, Q0W5T3, Q110Z9, Q11WK1, Q12EF7, Q12KC8, Q12WM7, Q13V40, Q14HS1, Q15YL3, Q17X40, Q18ES5, Q1AUK6, Q1BE73, Q1BYB1, Q1C3X1, Q1CLU7, Q1CUJ7, Q1DFA7, Q1GXW0, Q1I4H5, Q1IHV2, Q1IWZ8, Q1LQM3, Q1MPW7, Q1QEE8, Q1QSB0, Q1RG34, Q21MJ0, Q21YQ0, Q24VD4, Q254Z2, Q2A2U8, Q2FFN1, Q2FG69, Q2FTK5, Q2FXR4, Q2G283, Q2IHP8, Q2JFQ1, Q2JMP7, Q2JS70, Q2KUS4, Q2LVW6, Q2NF42, Q2NVQ0, Q2P6E4, Q2RJ31, Q2S1S3, Q2S8X9, Q2SYB2, Q2Y5Q5, Q2YTC2, Q2YU22, Q30RJ4, Q30WH2, Q31C50, Q31FB1, Q31QJ2, Q325Y5, Q32JV4, Q39566, Q39IQ4, Q39QA6, Q3A7W5, Q3ACS9, Q3ALU9, Q3AP59, Q3AWP4, Q3B1A1, Q3BPP7, Q3IHS1, Q3IT20, Q3JDH7, Q3JPN1, Q3K6C0, Q3KMF2, Q3M3B9, Q3SFU6, Q3Z5K3, Q40147, Q42522, Q46CH1, Q46GT9, Q46XZ7, Q478V1, Q47LB3, Q47V96, Q48DP1, Q49Y99, Q49YQ9, Q4FVC4, Q4JAM7, Q4JSY1, Q4K5I6, Q4L706, Q4L7G9, Q4UYA7, Q4ZNA0, Q55665, Q57T53, Q58020, Q5E2W6, Q5FAH9, Q5H3I5, Q5HER0, Q5HFA5, Q5HN71, Q5HNN6, Q5HUU3, Q5KWK9, Q5L2S4, Q5L5P0, Q5NGB9, Q5P6J6, Q5PD43, Q5QVP9, Q5SJS4, Q5UZ90, Q5WEP8, Q5WIE9, Q5WWG1, Q5X529, Q5YP79, Q5ZVA6, Q60CV0, Q62HV8, Q633Y3, Q63GB4, Q63RP8, Q65GK4, Q65M82, Q65U43, Q66EF1, Q67KM3, Q6AB08, Q6AHE8, Q6AQ32, Q6D1Z0, Q6FCY1, Q6G870, Q6G8Q8, Q6GFJ3, Q6GG38, Q6HD65, Q6HNS8, Q6L2G9, Q6LUS3, Q6M0P5, Q6MAC7, Q6MHT9, Q6NJJ2, Q6YZE2, Q71YY2, Q71ZB5, Q725I1, Q72K83, Q72ZW5, Q73DX4, Q73SQ3, Q74GA9, Q75G04, Q7A4T5, Q7M847, Q7MHY9, Q7N845, Q7NPI4, Q7P1Z5, Q7U598, Q7UPM9, Q7V2J3, Q7V677, Q7VDA1, Q7VHK3, Q7W050, Q7W3U1, Q7WF71, Q817R3, Q81I85, Q81LD0, Q81YV0, Q821C1, Q822Q0, Q82E21, Q82UQ8, Q83AK3, Q83FJ5, Q83H98, Q85WB7, Q87BW3, Q87LY3, Q87VY5, Q88DP0, Q897K4, Q8CNZ1, Q8CRW7, Q8CV56, Q8CZE3, Q8D3C8, Q8DBX8, Q8DLK8, Q8EHC8, Q8EY44, Q8FL16, Q8FSD4, Q8KAQ7, Q8NT73, Q8NVU6, Q8NW75, Q8P5R4, Q8PH40, Q8PW58, Q8RFY7, Q8TT57, Q8TYL6, Q8X4V5, Q8Y1M4, Q8Y6J9, Q8Y6X8, Q8YS26, Q8Z9B4, Q8ZBL9, Q8ZYW1, Q92AX5, Q92BG1, Q976H2, Q97B25, Q97MU2, Q980U5, Q99T15, Q9CNG9, Q9F2S0, Q9HKM6, Q9HMY8, Q9JRW9, Q9JW10, Q9JXW0, Q9K8G3, Q9KEB0, Q9KU97, Q9PB43, Q9PKI3, Q9PP70, Q9RWW0, Q9Y9I9, Q9ZMD0] | 5.4.3.8_1: [] | "
$ws.Range("H3").Value = "RXN0-1382: [A1ADQ1, A4YXN2, A5EGD7, A7ZPI2, A8A2M8, B1IX88, B1LMH0, B1X9P6, B2TWX3, B3QBS6, B5YYX4, B6I6S5, B6JE29, B7LBS7, B7M6P3, B7MH34, B7MY33, B7N5X4, B7NPQ8, B7UG84, C4ZVR1, O06644, O87838, P69902, P69903, Q07Q82, Q0TF87, Q139H7, Q1R8Z2, Q217M3, Q2IUI7, Q31Y97, Q32DG9, Q3YZF6, Q6N8F8, Q82M40, Q89QH2, Q8FFE8, Q8XBR7] | R07290: [] | 2.8.3.16_1: [] | "
$ws.Range("F4").Value = "RXN-12074: [Q88JX9] | R00008_reverse: [] | 4.1.3.17_7: [] | 4.1.3.17-RXN: [] | "
$ws.Range("G4").Value = "R02272: [] | GSAAMINOTRANS-RXN: [A0AJ02, A0AJD9, A0B885, A0JYM2, A0K5D7, A0KNY9, A0KZS6, A0L3S9, A0LIP3, A0LRE5, A0LXY3, A0PLS0, A0Q2A8, A0Q628, A0QLG2, A0QR33, A0R9G6, A0RJ78, A0RQ51, A0RXB3, A1A7K0, A1ASE9, A1AVQ1, A1BJG8, A1JJQ1, A1KB59, A1KG00, A1KS34, A1R898, A1RMG7, A1RQR5, A1S3U0, A1SE06, A1ST92, A1T3F2, A1TKK8, A1TYG5, A1UAR1, A1V1L0, A1V9X6, A1VSM6, A1VZJ6, A1WBR8, A1WMA0, A1WYL3, A2BL27, A2BPW6, A2BVE5, A2C0U2, A2C7I7, A2S9D9, A2SKQ7, A2SQU3, A3CL81, A3CU65, A3D1R5, A3DIF3, A3M7I8, A3MMQ8, A3MWW7, A3N2K3, A3NCF3, A3NY82, A3PBK8, A3PUB7, A3QBN5, A4FPX3, A4FZY1, A4G8N0, A4IKL1, A4IRG2, A4IXN8, A4J6H0, A4JC86, A4QB78, A4SGT2, A4SJ79, A4SVE6, A4T358, A4TPW6, A4VQY0, A4W6Q1, A4WN33, A4XC73, A4XK09, A4XYV6, A4Y4G4, A4YD91, A5CNI7, A5CXH8, A5D3M2, A5EXX1, A5F945, A5FIM3, A5G9C0, A5GMT2, A5GUJ2, A5IC29, A5ITJ2, A5IU30, A5N5Y1, A5TZQ4, A5UML0, A5UU40, A5VM65, A5W9H0, A5WC94, A6GYW4, A6LKX4, A6LSY8, A6Q2X6, A6Q9P9, A6QHK1, A6QI46, A6T2A2, A6T4V8, A6TJD8, A6U2D5, A6U2W8, A6URL2, A6UWB3, A6V0D2, A6VIL1, A6VQ66, A6VUX8, A6W5M6, A6WKL6, A7FM09, A7GL00, A7GTE2, A7GXK4, A7H3N1, A7HIX3, A7I252, A7I7P2, A7MGR5, A7MUU9, A7NCX8, A7NKV1, A7X388, A7X3X9, A7Z2N8, A7Z7A2, A7ZCH5, A7ZHP6, A7ZWA2, A8AAB3, A8ALD5, A8ETJ2, A8FB92, A8FFU9, A8FLR2, A8FS76, A8G3J9, A8G9U7, A8H164, A8LCD0, A8M6S5, A8M919, A8MGY8, A8YY36, A8Z2I8, A8ZZU9, A9A1A0, A9A875, A9AEU0, A9B550, A9BEA5, A9BY74, A9GDD3, A9I1R3, A9KBA0, A9KP86, A9KSC7, A9L5J5, A9M1G6, A9MPK7, A9N0P9, A9NAX3, A9R1E6, A9VIT5, A9VSS7, A9WIS7, A9WPI0, B0B9W0, B0BBJ0, B0BRF0, B0CC57, B0JPW6, B0KJV9, B0R7L9, B0RBL3, B0RP74, B0SF64, B0SNK3, B0TFV0, B0TIQ0, B0U0T6, B0U328, B0UST3, B0VCN6, B0VU92, B1HUT1, B1HVD3, B1I4L1, B1IQI6, B1J1Y6, B1JK22, B1JX81, B1KHF6, B1LGV7, B1MHH1, B1VEK3, B1VTD1, B1X023, B1XD24, B1XIT5, B1XT79, B1Y8G2, B1YBL3, B1YJU8, B1YK58, B1YUJ7, B1ZVF7, B2A1H4, B2FT35, B2G9H2, B2GL27, B2HRL5, B2HW19, B2I6C2, B2J7M9, B2JED5, B2K547, B2SGG7, B2SKS0, B2T6C2, B2TPD4, B2U2Z9, B2U6Q7, B2UNE3, B2USD7, B2UX41, B2V5U0, B2VE25, B3DY05, B3E9X0, B3EI07, B3EKJ7, B3GYN8, B3PL68, B3QRD2, B3QSA6, B3R312, B4EBT3, B4EUE1, B4RNX2, B4S0D0, B4S3Q6, B4SGW1, B4SI23, B4SUY4, B4TK30, B4TXQ6, B4U8M6, B4UKE6, B5BL82, B5EFG4, B5ERN9, B5F8R5, B5FAL8, B5FJ01, B5R3G6, B5RHE0, B5Y1L5, B5YHH6, B5Z0D4, B5ZA72, B6EL04, B6HZD0, B6J3H2, B6J495, B6JKN5, B7GH35, B7GIK0, B7GYZ5, B7H9S6, B7HE93, B7HQM1, B7HU66, B7I4N4, B7IHH3, B7IIX2, B7IW22, B7JBI1, B7JNG6, B7JQ54, B7K2I1, B7KA18, B7LGL6, B7LWB5, B7M195, B7MBD7, B7MP16, B7N823, B7NIB7, B7UIK1, B7V9D9, B7VJJ5, B8CJG4, B8DFN7, B8DHJ5, B8DN02, B8EBU3, B8F7H5, B8FBA1, B8G2L6, B8G822, B8GJJ2, B8GLE4, B8HBR2, B8HYK1, B8I0R1, B8J3A3, B8JCU6, B8ZT67, B9DMX3, B9DND8, B9DZG0, B9E740, B9E7Y3, B9IZ36, B9J3N2, B9KFS8, B9L0Q2, B9LA83, B9LKS0, B9LTZ9, B9LZL6, B9MF68, B9MRJ7, C0Q5R5, C0QCR9, C0QQ98, C0Z671, C0ZAH4, C0ZUX9, C1A7K7, C1AKK4, C1AZ30, C1CWI8, C1D9B6, C1DMY5, C1DUY4, C1ETQ7, C1EWG6, C1F911, C1KVK0, C1KVY4, C3K2K3, C3L6Z2, C3LHA0, C3LSN1, C3MJ22, C3MYE0, C3MZR9, C3N842, C3NF69, C3P9E6, C3PCZ3, C3PKI4, C4K6Y5, C4KJ00, C4L367, C4L4J9, C4LAK4, C4ZRP7, C5BAP9, C6AR33, C6DC32, E1W874, O26330, O29027, O66998, O74038, O74061, O84212, P0C1P8, P0C2D9, P0CL07, P18492, P21267, P23893, P24630, P30949, P31593, P45621, P46716, P48247, P56115, P63507, P63508, P71084, P99096, P9WMN8, P9WMN9, Q01YQ2, Q02SE5, Q04NV4, Q04X52, Q06741, Q06774, Q07YU5, Q0AAH7, Q0AFV1, Q0AZ36, Q0BHJ2, Q0BLC7, Q0HG53, Q0HSE6, Q0I3R7, Q0I8G1, Q0KDN9, Q0RS03, Q0SF42, Q0ST19, Q0T852, Q0TLH7, Q0TQG7, Q0VSP6, Q0W5T3, Q110Z9, Q11WK1, Q12EF7, Q12KC8, Q12WM7, Q13V40, Q14HS1, Q15YL3, Q17X40, Q18ES5, Q1AUK6, Q1BE73, Q1BYB1, Q1C3X1, Q1CLU7, Q1CUJ7, Q1DFA7, Q1GXW0, Q1I4H5, Q1IHV2, Q1IWZ8, Q1LQM3, Q1MPW7, Q1QEE8, Q1QSB0, Q1RG34, Q21MJ0, Q21YQ0, Q24VD4, Q254Z2, Q2A2U8, Q2FFN1, Q2FG69, Q2FTK5, Q2FXR4, Q2G283, Q2IHP8, Q2JFQ1, Q2JMP7, Q2JS70, Q2KUS4, Q2LVW6, Q2NF42, Q2NVQ0, Q2P6E4, Q2RJ31, Q2S1S3, Q2S8X9, Q2SYB2, Q2Y5Q5, Q2YTC2, Q2YU22, Q30RJ4, Q30WH2, Q31C50, Q31FB1, Q31QJ2, Q325Y5, Q32JV4, Q39566, Q39IQ4, Q39QA6, Q3A7W5, Q3ACS9, Q3ALU9, Q3AP59, Q3AWP4, Q3B1A1, Q3BPP7, Q3IHS1, Q3IT20, Q3JDH7, Q3JPN1, Q3K6C0, Q3KMF2, Q3M3B9, Q3SFU6, Q3Z5K3, Q40147, Q42522, Q46CH1, Q46GT9, Q46XZ7, Q478V1, Q47LB3, Q47V96, Q48DP1, Q49Y99, Q49YQ9, Q4FVC4, Q4JAM7, Q4JSY1, Q4K5I6, Q4L706, Q4L7G9, Q4UYA7, Q4ZNA0, Q55665, Q57T53, Q58020, Q5E2W6, Q5FAH9, Q5H3I5, Q5HER0, Q5HFA5, Q5HN71, Q5HNN6, Q5HUU3, Q5KWK9, Q5L2S4, Q5L5P0, Q5NGB9, Q5P6J6, Q5PD43, Q5QVP9, Q5SJS4, Q5UZ90, Q5WEP8, Q5WIE9, Q5WWG1, Q5X529, Q5YP79, Q5ZVA6, Q60CV0, Q62HV8, Q633Y3, Q63GB4, Q63RP8, Q65GK4, Q65M82, Q65U43, Q66EF1, Q67KM3, Q6AB08, Q6AHE8, Q6AQ32, Q6D1Z0, Q6FCY1, Q6G870, Q6G8Q8, Q6GFJ3, Q6GG38, Q6HD65, Q6HNS8, Q6L2G9, Q6LUS3, Q6M0P5, Q6MAC7, Q6MHT9, Q6NJJ2, Q6YZE2, Q71YY2, Q71ZB5, Q725I1, Q72K83, Q72ZW5, Q73DX4, Q73SQ3, Q74GA9, Q75G04, Q7A4T5, Q7M847, Q7MHY9, Q7N845, Q7NPI4, Q7P1Z5, Q7U598, Q7UPM9, Q7V2J3, Q7V677, Q7VDA1, Q7VHK3, Q7W050, Q7W3U1, Q7WF71, Q817R3, Q81I85, Q81LD0, Q81YV0, Q821C1, Q822Q0, Q82E21, Q82UQ8, Q83AK3, Q83FJ5, Q83H98, Q85WB7, Q87BW3, Q87LY3, Q87VY5, Q88DP0, Q897K4, Q8CNZ1, Q8CRW7, Q8CV56, Q8CZE3, Q8D3C8, Q8DBX8, Q8DLK8, Q8EHC8, Q8EY44, Q8FL16, Q8FSD4, Q8KAQ7, Q8NT73, Q8NVU6, Q8NW75, Q8P5R4, Q8PH40, Q8PW58, Q8RFY7, Q8TT57, Q8TYL6, Q8X4V5, Q8Y1M4, Q8Y6J9, Q8Y6X8, Q8YS26, Q8Z9B4, Q8ZBL9, Q8ZYW1, Q92AX5, Q92BG1, Q976H2, Q97B25, Q97MU2, Q980U5, Q99T15, Q9CNG9, Q9F2S0, Q9HKM6, Q9HMY8, Q9JRW9, Q9JW10, Q9JXW0, Q9K8G3, Q9KEB0, Q9KU97, Q9PB43, Q9PKI3, Q9PP70, Q9RWW0, Q9Y9I9, Q9ZMD0] | 5.4.3.8_1: [] | "
$ws.Range("H4").Value = "4.1.1.4_1: [] | "

